$wb = $excel.ActiveWorkbook

# --- UITestingButton sheet: collapse to a single column of names ---
$ws2 = $wb.Worksheets.Item("UITestingButton")
$ws2.Range("A1").Value = "Sharmin"
$ws2.Range("A2").Value = "Mehnaz"
$ws2.Range("A3").Value = "Nusrat"
$ws2.Range("B1:B3").ClearContents()
[void]$ws2.Range("N10").Select()

# --- New sheet: OrangeHRMSearches ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "OrangeHRMSearches"
$ws3.Range("A1").Value = "Search"
$ws3.Range("A2").Value = "Leave"
$ws3.Range("A3").Value = "Performance"
[void]$ws3.Activate()
[void]$ws3.Range("A3").Select()
